# Apply scheduled runner updates to Sheets/Famfrit_Profits.xlsx
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) for the
# affected leve rows across all job sheets, as produced by the pricing refresh job.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33 (Leve Item ID 5512)
$ws.Cells.Item(33, 8).Value = 13065.883
$ws.Cells.Item(33, 9).Value = 19710.818
$ws.Cells.Item(33, 10).Value = 883.5
$ws.Cells.Item(33, 11).Value = 19710.818
$ws.Cells.Item(33, 12).Value = 883.5
$ws.Cells.Item(33, 13).Value = -19481.818
$ws.Cells.Item(33, 14).Value = -1341.5
# Row 41 (Leve Item ID 5478)
$ws.Cells.Item(41, 8).Value = 498.33334
$ws.Cells.Item(41, 9).Value = 500
$ws.Cells.Item(41, 10).Value = 495
$ws.Cells.Item(41, 11).Value = 500
$ws.Cells.Item(41, 12).Value = 495
$ws.Cells.Item(41, 13).Value = -60
$ws.Cells.Item(41, 14).Value = -1375
# Row 76 (Leve Item ID 12602)
$ws.Cells.Item(76, 8).Value = 9073.071
$ws.Cells.Item(76, 9).Value = 10007.143
$ws.Cells.Item(76, 10).Value = 8139
$ws.Cells.Item(76, 11).Value = 10007.143
$ws.Cells.Item(76, 12).Value = 8139
$ws.Cells.Item(76, 13).Value = -9692.143
# Row 79 (Leve Item ID 12602)
$ws.Cells.Item(79, 8).Value = 9073.071
$ws.Cells.Item(79, 9).Value = 10007.143
$ws.Cells.Item(79, 10).Value = 8139
$ws.Cells.Item(79, 11).Value = 10007.143
$ws.Cells.Item(79, 12).Value = 8139
$ws.Cells.Item(79, 13).Value = -8915.143
# Row 97 (Leve Item ID 19885)
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 14).ClearContents()
# Row 112 (Leve Item ID 27960)
$ws.Cells.Item(112, 8).Value = 14708448
$ws.Cells.Item(112, 9).Value = 5187.5
$ws.Cells.Item(112, 10).Value = 15627402
$ws.Cells.Item(112, 11).Value = 15562.5
$ws.Cells.Item(112, 12).Value = 46882206
$ws.Cells.Item(112, 13).Value = -14454.5
$ws.Cells.Item(112, 14).Value = -46884422
# Row 125 (Leve Item ID 36228)
$ws.Cells.Item(125, 8).Value = 8762
$ws.Cells.Item(125, 9).Value = 11000
$ws.Cells.Item(125, 10).Value = 8016
$ws.Cells.Item(125, 11).Value = 99000
$ws.Cells.Item(125, 12).Value = 72144
$ws.Cells.Item(125, 13).Value = -96540
$ws.Cells.Item(125, 14).Value = -77064
# Row 132 (Leve Item ID 44049)
$ws.Cells.Item(132, 8).Value = 5221.074
$ws.Cells.Item(132, 9).Value = 6013.727
$ws.Cells.Item(132, 10).Value = 1733.4
$ws.Cells.Item(132, 11).Value = 18041.181
$ws.Cells.Item(132, 12).Value = 5200.200000000001
$ws.Cells.Item(132, 13).Value = -15511.181
$ws.Cells.Item(132, 14).Value = -10260.2
# Row 135 (Leve Item ID 44047)
$ws.Cells.Item(135, 8).Value = 450.8846
$ws.Cells.Item(135, 9).Value = 421.79166
$ws.Cells.Item(135, 10).Value = 800
$ws.Cells.Item(135, 11).Value = 3796.12494
$ws.Cells.Item(135, 12).Value = 7200
$ws.Cells.Item(135, 13).Value = -1261.12494

$ws = $wb.Worksheets.Item("ARM")
# Row 45 (Leve Item ID 27714)
$ws.Cells.Item(45, 8).Value = 3733.3333
$ws.Cells.Item(45, 9).Value = 2980
$ws.Cells.Item(45, 10).Value = 7500
$ws.Cells.Item(45, 11).Value = 2980
$ws.Cells.Item(45, 12).Value = 7500
$ws.Cells.Item(45, 13).Value = -2603
# Row 57 (Leve Item ID 39767)
$ws.Cells.Item(57, 8).Value = 9999
$ws.Cells.Item(57, 9).Value = 9999
$ws.Cells.Item(57, 10).Value = 0
$ws.Cells.Item(57, 11).Value = 9999
$ws.Cells.Item(57, 12).Value = 0
$ws.Cells.Item(57, 13).Value = -9515
# Row 61 (Leve Item ID 43999)
$ws.Cells.Item(61, 8).Value = 2551.3076
$ws.Cells.Item(61, 9).Value = 2112.348
$ws.Cells.Item(61, 10).Value = 5916.6665
$ws.Cells.Item(61, 11).Value = 2112.348
$ws.Cells.Item(61, 12).Value = 5916.6665
$ws.Cells.Item(61, 13).Value = -1900.348
$ws.Cells.Item(61, 14).Value = -6340.6665
# Row 74 (Leve Item ID 44000)
$ws.Cells.Item(74, 8).Value = 4934.7715
$ws.Cells.Item(74, 9).Value = 4749
$ws.Cells.Item(74, 10).Value = 8000
$ws.Cells.Item(74, 11).Value = 4749
$ws.Cells.Item(74, 12).Value = 8000
$ws.Cells.Item(74, 13).Value = -3875
# Row 77 (Leve Item ID 44000)
$ws.Cells.Item(77, 8).Value = 4934.7715
$ws.Cells.Item(77, 9).Value = 4749
$ws.Cells.Item(77, 10).Value = 8000
$ws.Cells.Item(77, 11).Value = 23745
$ws.Cells.Item(77, 12).Value = 40000
$ws.Cells.Item(77, 13).Value = -19377
# Row 97 (Leve Item ID 19941)
$ws.Cells.Item(97, 8).Value = 1568.7273
$ws.Cells.Item(97, 9).Value = 1450.6316
$ws.Cells.Item(97, 10).Value = 2316.6667
$ws.Cells.Item(97, 11).Value = 1450.6316
$ws.Cells.Item(97, 12).Value = 2316.6667
$ws.Cells.Item(97, 13).Value = -954.6315999999999
# Row 102 (Leve Item ID 19945)
$ws.Cells.Item(102, 8).Value = 3966
$ws.Cells.Item(102, 9).Value = 3484.5
$ws.Cells.Item(102, 10).Value = 5699.4
$ws.Cells.Item(102, 11).Value = 3484.5
$ws.Cells.Item(102, 12).Value = 5699.4
$ws.Cells.Item(102, 13).Value = -1862.5
# Row 132 (Leve Item ID 43997)
$ws.Cells.Item(132, 8).Value = 38632.207
$ws.Cells.Item(132, 9).Value = 2577.7058
$ws.Cells.Item(132, 10).Value = 283802.8
$ws.Cells.Item(132, 11).Value = 7733.117400000001
$ws.Cells.Item(132, 12).Value = 851408.3999999999
$ws.Cells.Item(132, 13).Value = -5203.117400000001
$ws.Cells.Item(132, 14).Value = -856468.3999999999
# Row 136 (Leve Item ID 43999)
$ws.Cells.Item(136, 8).Value = 2551.3076
$ws.Cells.Item(136, 9).Value = 2112.348
$ws.Cells.Item(136, 10).Value = 5916.6665
$ws.Cells.Item(136, 11).Value = 6337.044
$ws.Cells.Item(136, 12).Value = 17749.9995
$ws.Cells.Item(136, 13).Value = -3787.044
$ws.Cells.Item(136, 14).Value = -22849.9995

$ws = $wb.Worksheets.Item("BSM")
# Row 64 (Leve Item ID 14184)
$ws.Cells.Item(64, 8).Value = 419
$ws.Cells.Item(64, 9).Value = 734.6667
$ws.Cells.Item(64, 10).Value = 103.333336
$ws.Cells.Item(64, 11).Value = 734.6667
$ws.Cells.Item(64, 12).Value = 103.333336
$ws.Cells.Item(64, 13).Value = -509.6667
$ws.Cells.Item(64, 14).Value = -553.333336
# Row 67 (Leve Item ID 14184)
$ws.Cells.Item(67, 8).Value = 419
$ws.Cells.Item(67, 9).Value = 734.6667
$ws.Cells.Item(67, 10).Value = 103.333336
$ws.Cells.Item(67, 11).Value = 734.6667
$ws.Cells.Item(67, 12).Value = 103.333336
$ws.Cells.Item(67, 13).Value = 45.33330000000001
$ws.Cells.Item(67, 14).Value = -1663.333336
# Row 82 (Leve Item ID 11877)
$ws.Cells.Item(82, 8).Value = 25460.625
$ws.Cells.Item(82, 9).Value = 14456
$ws.Cells.Item(82, 10).Value = 58474.5
$ws.Cells.Item(82, 11).Value = 14456
$ws.Cells.Item(82, 12).Value = 58474.5
$ws.Cells.Item(82, 13).Value = -14073
$ws.Cells.Item(82, 14).Value = -59240.5
# Row 85 (Leve Item ID 11877)
$ws.Cells.Item(85, 8).Value = 25460.625
$ws.Cells.Item(85, 9).Value = 14456
$ws.Cells.Item(85, 10).Value = 58474.5
$ws.Cells.Item(85, 11).Value = 14456
$ws.Cells.Item(85, 12).Value = 58474.5
$ws.Cells.Item(85, 13).Value = -13130
$ws.Cells.Item(85, 14).Value = -61126.5
# Row 94 (Leve Item ID 19939)
$ws.Cells.Item(94, 8).Value = 1402
$ws.Cells.Item(94, 9).Value = 842.6111
$ws.Cells.Item(94, 10).Value = 2317.3635
$ws.Cells.Item(94, 11).Value = 842.6111
$ws.Cells.Item(94, 12).Value = 2317.3635
$ws.Cells.Item(94, 13).Value = -391.6111

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Cells.Item(31, 8).Value = 4206.976
$ws.Cells.Item(31, 9).Value = 2249.5293
$ws.Cells.Item(31, 10).Value = 5538.04
$ws.Cells.Item(31, 11).Value = 2249.5293
$ws.Cells.Item(31, 12).Value = 5538.04
$ws.Cells.Item(31, 13).Value = -1954.5293
# Row 34 (Leve Item ID 44023)
$ws.Cells.Item(34, 8).Value = 4206.976
$ws.Cells.Item(34, 9).Value = 2249.5293
$ws.Cells.Item(34, 10).Value = 5538.04
$ws.Cells.Item(34, 11).Value = 2249.5293
$ws.Cells.Item(34, 12).Value = 5538.04
$ws.Cells.Item(34, 13).Value = -2047.5293
# Row 132 (Leve Item ID 44019)
$ws.Cells.Item(132, 8).Value = 4924.76
$ws.Cells.Item(132, 9).Value = 4847.8823
$ws.Cells.Item(132, 10).Value = 5088.125
$ws.Cells.Item(132, 11).Value = 14543.6469
$ws.Cells.Item(132, 12).Value = 15264.375
$ws.Cells.Item(132, 13).Value = -12013.6469

$ws = $wb.Worksheets.Item("CUL")
# Row 131 (Leve Item ID 36060)
$ws.Cells.Item(131, 8).Value = 1271.8572
$ws.Cells.Item(131, 9).Value = 1039
$ws.Cells.Item(131, 10).Value = 1446.5
$ws.Cells.Item(131, 11).Value = 3117
$ws.Cells.Item(131, 12).Value = 4339.5
$ws.Cells.Item(131, 13).Value = 1923
$ws.Cells.Item(131, 14).Value = -14419.5
# Row 132 (Leve Item ID 43972)
$ws.Cells.Item(132, 8).Value = 3091.8462
$ws.Cells.Item(132, 9).Value = 1533
$ws.Cells.Item(132, 10).Value = 3559.5
$ws.Cells.Item(132, 11).Value = 13797
$ws.Cells.Item(132, 12).Value = 32035.5
$ws.Cells.Item(132, 13).Value = -11267
$ws.Cells.Item(132, 14).Value = -37095.5
# Row 133 (Leve Item ID 44073)
$ws.Cells.Item(133, 8).Value = 1752.6
$ws.Cells.Item(133, 9).Value = 910
$ws.Cells.Item(133, 10).Value = 3016.5
$ws.Cells.Item(133, 11).Value = 2730
$ws.Cells.Item(133, 12).Value = 9049.5
$ws.Cells.Item(133, 13).Value = 2330

$ws = $wb.Worksheets.Item("GSM")
# Row 122 (Leve Item ID 36182)
$ws.Cells.Item(122, 8).Value = 2183
$ws.Cells.Item(122, 9).Value = 1519.6
$ws.Cells.Item(122, 10).Value = 5500
$ws.Cells.Item(122, 11).Value = 4558.799999999999
$ws.Cells.Item(122, 12).Value = 16500
$ws.Cells.Item(122, 13).Value = -2108.799999999999
$ws.Cells.Item(122, 14).Value = -21400

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Cells.Item(7, 8).Value = 5158.5557
$ws.Cells.Item(7, 9).Value = 4078.2727
$ws.Cells.Item(7, 10).Value = 6856.143
$ws.Cells.Item(7, 11).Value = 4078.2727
$ws.Cells.Item(7, 12).Value = 6856.143
$ws.Cells.Item(7, 13).Value = -3966.2727
# Row 40 (Leve Item ID 36248)
$ws.Cells.Item(40, 8).Value = 4062.7334
$ws.Cells.Item(40, 9).Value = 3367.375
$ws.Cells.Item(40, 10).Value = 4857.4287
$ws.Cells.Item(40, 11).Value = 3367.375
$ws.Cells.Item(40, 12).Value = 4857.4287
$ws.Cells.Item(40, 13).Value = -3231.375
# Row 122 (Leve Item ID 36247)
$ws.Cells.Item(122, 8).Value = 4366.1816
$ws.Cells.Item(122, 9).Value = 4091.125
$ws.Cells.Item(122, 10).Value = 5099.6665
$ws.Cells.Item(122, 11).Value = 12273.375
$ws.Cells.Item(122, 12).Value = 15298.9995
$ws.Cells.Item(122, 13).Value = -9823.375
# Row 126 (Leve Item ID 36249)
$ws.Cells.Item(126, 8).Value = 5158.5557
$ws.Cells.Item(126, 9).Value = 4078.2727
$ws.Cells.Item(126, 10).Value = 6856.143
$ws.Cells.Item(126, 11).Value = 12234.8181
$ws.Cells.Item(126, 12).Value = 20568.429
$ws.Cells.Item(126, 13).Value = -9764.8181
# Row 136 (Leve Item ID 44060)
$ws.Cells.Item(136, 8).Value = 4584.25
$ws.Cells.Item(136, 9).Value = 4692.8184
$ws.Cells.Item(136, 10).Value = 3390
$ws.Cells.Item(136, 11).Value = 14078.4552
$ws.Cells.Item(136, 12).Value = 10170
$ws.Cells.Item(136, 13).Value = -11528.4552
$ws.Cells.Item(136, 14).Value = -15270

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (Leve Item ID 44029)
$ws.Cells.Item(132, 8).Value = 1012.4103
$ws.Cells.Item(132, 9).Value = 864.19354
$ws.Cells.Item(132, 10).Value = 1586.75
$ws.Cells.Item(132, 11).Value = 2592.58062
$ws.Cells.Item(132, 12).Value = 4760.25
$ws.Cells.Item(132, 13).Value = -62.58061999999973

Write-Host "Famfrit_Profits: updated 232 cells across 8 sheets"
